# Apply the "Código y actualización tabla" update to the sensor table.
# Row 5 of the table ("Hoja1"!A5:E5) is updated: the LDR / Sensor Pasivo
# light sensor entry is replaced with a KY-018 / Sensor Activo entry, and
# the additional-components cell now reads "No requiere de otros ..." like
# the other active-sensor rows (so it also needs the wrap-text style used
# by column E elsewhere in the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5: sensor name
$ws.Range("A5").Value = "KY-018"

# B5 (Variable) stays "Luminosidad " - unchanged

# C5: sensor type
$ws.Range("C5").Value = "Sensor Activo"

# D5: ports
$ws.Range("D5").Value = "D22`n3V3`nGND"

# E5: additional components text + matching wrap style used by other
# "No requiere de otros..." cells (e.g. E4)
$ws.Range("E5").Value = "No requiere de otros `ncomponentes "
$ws.Range("E5").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E5").VerticalAlignment = -4108    # xlCenter
$ws.Range("E5").WrapText = $true

# Update the active selection to match the saved workbook state (D6).
$ws.Range("D6").Select()
